$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 0.4510973247164258
$ws.Range("J2").Value = 0.4510973247164258
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4977036666666666
$ws.Range("N2").Value = 1.493111
$ws.Range("Q2").Value = 0.1012088701227778
$ws.Range("R2").Value = 0.9108798311049999
$ws.Range("S2").Value = 0.4510973247164258
$ws.Range("T2").Value = 0.4510973247164258

# Row 3 updates
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2474416666666667
$ws.Range("H3").Value = 0.742325
$ws.Range("I3").Value = 0.5489026752835741
$ws.Range("J3").Value = 0.5489026752835741
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.4977036666666666
$ws.Range("N3").Value = 1.493111
$ws.Range("Q3").Value = 0.1231526247861111
$ws.Range("R3").Value = 1.108373623075
$ws.Range("S3").Value = 0.5489026752835741
$ws.Range("T3").Value = 0.5489026752835741
